$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Duplicate the existing "2022-Q2" sheet, dropping the copy right
#    after it. The ORIGINAL sheet (which keeps its original rId/file)
#    becomes the new "2022-Q3" sheet; the COPY becomes the new
#    "2022-Q2" sheet (identical data to the original, just re-homed
#    on a new sheet part) - mirroring how Excel assigns sheetIds when
#    a sheet is duplicated.
# ------------------------------------------------------------------
$qSheet = $wb.Worksheets.Item("2022-Q2")
$qSheet.Copy($null, $qSheet)

$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

$q2copy = $wb.Worksheets.Item(3)
$q2copy.Name = "2022-Q2"

# ------------------------------------------------------------------
# 2. Refresh the fund holdings figures on the (now) "2022-Q3" sheet
#    with the new quarter's numbers. The source data keeps these as
#    text (same as the original sheet), so format the range as text
#    before writing to stop Excel from auto-coercing to numbers.
# ------------------------------------------------------------------
$q3.Range("D2:G4").NumberFormat = "@"
$q3.Range("D2:D4").Value = "0.64"
$q3.Range("E2:E4").Value = "80.13"
$q3.Range("F2:F4").Value = "2.06"
$q3.Range("G2:G4").Value = "0.0132"

# ------------------------------------------------------------------
# 3. Update the "总计" (summary) sheet: the top row now reflects the
#    new 2022-Q3 totals, and a new row is appended for the prior
#    2022-Q2 totals (pushed down from row 2 to row 3).
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Range("B2").Value = "2022-Q3"
$summary.Range("D2").Value = 0.04

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 3
$summary.Range("D3").Value = 0.03

# Match row 3's "A" column formatting (bold, centered, bordered) to
# row 2's, same as the rest of that styled column.
$summary.Range("A2").Copy()
$summary.Range("A3").PasteSpecial(-4122)
